$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text, not auto-converted to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.401.39'
$ws.Range("D3").Value = '1.869.53'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '243.68'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").Value = '0.7040'
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '0.07920'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("D9").Value = '0.3134'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("D11").Value = '0.07839'
$ws.Range("E11").Value = '  -4.59%  '
$ws.Range("D12").Value = '1.910.58'
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = '93.83'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").Value = '5.170'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '0.7013'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").Value = '6.519'
$ws.Range("E16").Value = '  +2.05%  '
$ws.Range("D17").Value = '0.000008399'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '29.490.60'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").Value = '252.43'
$ws.Range("E19").Value = '  +3.76%  '
$ws.Range("D20").Value = '2.146.00'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '7.672'
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -3.71%  '
$ws.Range("D26").Value = '9.011'
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = '161.67'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").Value = '4.317'
$ws.Range("D31").Value = '4.259'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("D33").Value = '0.05270'
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("D34").Value = '1.898'
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("D35").Value = '1.182'
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("D36").Value = '0.7518'
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").Value = '0.01879'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").Value = '1.274.18'
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").Value = '2.770'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = '0.8915'
$ws.Range("E41").Value = '  -2.01%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.045'
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '109.36'
$ws.Range("E43").Value = '  -3.38%  '
$ws.Range("D44").Value = '70.99'
$ws.Range("E44").Value = '  -4.44%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("E46").Value = '  -4.78%  '
$ws.Range("D47").Value = '2.035.58'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").Value = '9.635'
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("D49").Value = '1.804'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").Value = '0.5183'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").Value = '0.4305'
$ws.Range("E51").Value = '  -0.90%  '
